# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 159
$wsOff.Range("C3").Value = 114
$wsOff.Range("D3").Value = 49
$wsOff.Range("E3").Value = 23

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 160
$wsDef.Range("C3").Value = 121
$wsDef.Range("D3").Value = 39
$wsDef.Range("E3").Value = 16
